# Updated cryptos list — apply Price (D) and Volume(1h) (E) changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "28.019.80";  E = "  +1.86%  " },
    @{ Row = 3;  D = "1.647.04";   E = "  +1.76%  " },
    @{ Row = 4;  D = $null;        E = "  -0.01%  " },
    @{ Row = 5;  D = "213.49";     E = "  +1.25%  " },
    @{ Row = 7;  D = $null;        E = "  -0.06%  " },
    @{ Row = 8;  D = "23.42";      E = "  +2.67%  " },
    @{ Row = 9;  D = $null;        E = "  +1.52%  " },
    @{ Row = 10; D = $null;        E = "  +0.17%  " },
    @{ Row = 11; D = "0.0873";     E = "  -1.30%  " },
    @{ Row = 12; D = "1.880.50";   E = "  +1.73%  " },
    @{ Row = 13; D = "1.645.22";   E = "  +1.89%  " },
    @{ Row = 14; D = $null;        E = "  +1.01%  " },
    @{ Row = 15; D = "0.564";      E = "  +2.66%  " },
    @{ Row = 16; D = "65.58";      E = "  +0.60%  " },
    @{ Row = 17; D = "27.993.82";  E = "  +1.83%  " },
    @{ Row = 18; D = "232.65";     E = "  +0.77%  " },
    @{ Row = 19; D = $null;        E = "  +2.38%  " },
    @{ Row = 20; D = "0.0₃0723";   E = "  +0.63%  " },
    @{ Row = 21; D = $null;        E = "  +0.00%  " },
    @{ Row = 22; D = $null;        E = "  +4.96%  " },
    @{ Row = 23; D = "4.40";       E = "  +2.65%  " },
    @{ Row = 24; D = "2.14";       E = "  +3.58%  " },
    @{ Row = 25; D = "152.11";     E = "  +0.84%  " },
    @{ Row = 26; D = "6.92";       E = "  +1.24%  " },
    @{ Row = 28; D = "0.111";      E = "  +0.04%  " },
    @{ Row = 29; D = $null;        E = "  +0.05%  " },
    @{ Row = 30; D = "1.20";       E = "  +1.66%  " },
    @{ Row = 31; D = $null;        E = "  +0.29%  " },
    @{ Row = 32; D = $null;        E = "  +2.71%  " },
    @{ Row = 33; D = "1.444.01";   E = "  -1.58%  " },
    @{ Row = 34; D = $null;        E = "  +0.15%  " },
    @{ Row = 35; D = $null;        E = "  +1.96%  " },
    @{ Row = 36; D = $null;        E = "  -0.56%  " },
    @{ Row = 37; D = "0.889";      E = "  +3.49%  " },
    @{ Row = 38; D = $null;        E = "  +1.10%  " },
    @{ Row = 39; D = "0.561";      E = "  +0.35%  " },
    @{ Row = 40; D = "0.920";      E = "  -3.20%  " },
    @{ Row = 41; D = "69.31";      E = "  +2.14%  " },
    @{ Row = 42; D = $null;        E = "  +3.85%  " },
    @{ Row = 43; D = $null;        E = "  -0.05%  " },
    @{ Row = 44; D = "2.47";       E = "  -1.54%  " },
    @{ Row = 45; D = $null;        E = "  +0.99%  " },
    @{ Row = 46; D = $null;        E = "  +2.84%  " },
    @{ Row = 47; D = $null;        E = "  +4.75%  " },
    @{ Row = 48; D = "1.788.75";   E = "  +1.58%  " },
    @{ Row = 49; D = "88.93";      E = "  +2.68%  " },
    @{ Row = 50; D = $null;        E = "  -0.36%  " },
    @{ Row = 51; D = $null;        E = "  +0.18%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Leading apostrophe forces text interpretation (no NumberFormat
        # change) so numeric-looking prices like "213.49" stay exact text
        # instead of being coerced into doubles.
        $ws.Range("D$($u.Row)").Value = "'" + $u.D
    }
    $ws.Range("E$($u.Row)").Value = $u.E
}
